# Update "想去人数" (want-to-go count) values across the four worksheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 356
$ws1.Range("F7").Value = 1244
$ws1.Range("F9").Value = 2386
$ws1.Range("F10").Value = 834
$ws1.Range("F11").Value = 18268
$ws1.Range("F16").Value = 299
$ws1.Range("F17").Value = 582
$ws1.Range("F19").Value = 180
$ws1.Range("F25").Value = 67

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 102
$ws2.Range("F14").Value = 63

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5825

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5825
$ws4.Range("F8").Value = 356
$ws4.Range("F12").Value = 1244
$ws4.Range("F17").Value = 2386
$ws4.Range("F18").Value = 834
$ws4.Range("F19").Value = 18268
$ws4.Range("F22").Value = 102
$ws4.Range("F23").Value = 102
$ws4.Range("F28").Value = 299
$ws4.Range("F29").Value = 582
$ws4.Range("F31").Value = 180
$ws4.Range("F38").Value = 63
$ws4.Range("F50").Value = 67
